$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 (recomputed values) and add new row 5 (A1:BA5 dimension)

# A5 is a brand new cell; give it the same header-column style as A2:A4 (bold/border/center)
# before writing its value, by copying formats only from A4.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = [double]"0"
$ws.Range("B2").Value = [double]"0.01380925462327384"
$ws.Range("C2").Value = [double]"0.01446566436407964"
$ws.Range("D2").Value = [double]"0.01407997425376136"
$ws.Range("E2").Value = [double]"0.01483870967741936"
$ws.Range("F2").Value = [double]"0.01432710882163554"
$ws.Range("G2").Value = [double]"0.01460030652577236"
$ws.Range("H2").Value = [double]"0.01285553591515346"
$ws.Range("I2").Value = [double]"0.01423027166882277"
$ws.Range("J2").Value = [double]"0.01411421155094095"
$ws.Range("K2").Value = [double]"0.0128174123337364"
$ws.Range("L2").Value = [double]"0.01453581753935111"
$ws.Range("M2").Value = [double]"0.01360160965794769"
$ws.Range("N2").Value = [double]"0.01409804549823774"
$ws.Range("O2").Value = [double]"0.01536085825747724"
$ws.Range("P2").Value = [double]"0.01417525773195876"
$ws.Range("Q2").Value = [double]"0.01326646173758292"
$ws.Range("R2").Value = [double]"0.01461105215900743"
$ws.Range("S2").Value = [double]"0.01343648208469055"
$ws.Range("T2").Value = [double]"0.01345038659793815"
$ws.Range("U2").Value = [double]"0.01389893146943038"
$ws.Range("V2").Value = [double]"0.01438025982514911"
$ws.Range("W2").Value = [double]"0.01419292759201347"
$ws.Range("X2").Value = [double]"0.01388776090311053"
$ws.Range("Y2").Value = [double]"0.01355739568111706"
$ws.Range("Z2").Value = [double]"0.01374929645412881"
$ws.Range("AA2").Value = [double]"0.01364190012180268"
$ws.Range("AB2").Value = [double]"0.01445434834979523"
$ws.Range("AC2").Value = [double]"0.012524557956778"
$ws.Range("AD2").Value = [double]"0.01305714515999033"
$ws.Range("AE2").Value = [double]"0.01422417296594327"
$ws.Range("AF2").Value = [double]"0.01440945519307051"
$ws.Range("AG2").Value = [double]"0.01382041542067405"
$ws.Range("AH2").Value = [double]"0.01441614608361365"
$ws.Range("AI2").Value = [double]"0.01414563574987944"
$ws.Range("AJ2").Value = [double]"0.01342551652062063"
$ws.Range("AK2").Value = [double]"0.0132013201320132"
$ws.Range("AL2").Value = [double]"0.01271082864825226"
$ws.Range("AM2").Value = [double]"0.01425793907971484"
$ws.Range("AN2").Value = [double]"0.01343570057581574"
$ws.Range("AO2").Value = [double]"0.01398884127112477"
$ws.Range("AP2").Value = [double]"0.01347840207859695"
$ws.Range("AQ2").Value = [double]"0.01426844014510278"
$ws.Range("AR2").Value = [double]"0.01268994503718073"
$ws.Range("AS2").Value = [double]"0.01323861801743623"
$ws.Range("AT2").Value = [double]"0.01422943966556797"
$ws.Range("AU2").Value = [double]"0.01362309438858255"
$ws.Range("AV2").Value = [double]"0.01407877604166667"
$ws.Range("AW2").Value = [double]"0.01278918214402085"
$ws.Range("AX2").Value = [double]"0.01395084759510098"
$ws.Range("AY2").Value = [double]"0.01309262421728877"
$ws.Range("AZ2").Value = [double]"0.01382988570966739"
$ws.Range("BA2").Value = [double]"0.0006177913744012809"

# Row 3
$ws.Range("A3").Value = [double]"1"
$ws.Range("B3").Value = [double]"0.28125"
$ws.Range("C3").Value = [double]"0.2908496732026144"
$ws.Range("D3").Value = [double]"0.2864157119476268"
$ws.Range("E3").Value = [double]"0.3026315789473684"
$ws.Range("F3").Value = [double]"0.2908496732026144"
$ws.Range("G3").Value = [double]"0.2909967845659164"
$ws.Range("H3").Value = [double]"0.264026402640264"
$ws.Range("I3").Value = [double]"0.2953020134228188"
$ws.Range("J3").Value = [double]"0.2949152542372881"
$ws.Range("K3").Value = [double]"0.2602291325695581"
$ws.Range("L3").Value = [double]"0.2914653784219002"
$ws.Range("M3").Value = [double]"0.2752442996742671"
$ws.Range("N3").Value = [double]"0.2894736842105263"
$ws.Range("O3").Value = [double]"0.2976377952755905"
$ws.Range("P3").Value = [double]"0.2899505766062603"
$ws.Range("Q3").Value = [double]"0.2822719449225473"
$ws.Range("R3").Value = [double]"0.2934426229508197"
$ws.Range("S3").Value = [double]"0.2731788079470199"
$ws.Range("T3").Value = [double]"0.2715447154471545"
$ws.Range("U3").Value = [double]"0.2826797385620915"
$ws.Range("V3").Value = [double]"0.2866449511400652"
$ws.Range("W3").Value = [double]"0.3020477815699659"
$ws.Range("X3").Value = [double]"0.285"
$ws.Range("Y3").Value = [double]"0.2792642140468227"
$ws.Range("Z3").Value = [double]"0.2762520193861066"
$ws.Range("AA3").Value = [double]"0.2736156351791531"
$ws.Range("AB3").Value = [double]"0.2884615384615384"
$ws.Range("AC3").Value = [double]"0.2562814070351759"
$ws.Range("AD3").Value = [double]"0.2709030100334448"
$ws.Range("AE3").Value = [double]"0.2859477124183006"
$ws.Range("AF3").Value = [double]"0.2918032786885246"
$ws.Range("AG3").Value = [double]"0.301056338028169"
$ws.Range("AH3").Value = [double]"0.2936378466557912"
$ws.Range("AI3").Value = [double]"0.2866449511400652"
$ws.Range("AJ3").Value = [double]"0.2715447154471545"
$ws.Range("AK3").Value = [double]"0.2719734660033167"
$ws.Range("AL3").Value = [double]"0.2557377049180328"
$ws.Range("AM3").Value = [double]"0.2938230383973289"
$ws.Range("AN3").Value = [double]"0.2740619902120718"
$ws.Range("AO3").Value = [double]"0.2922297297297297"
$ws.Range("AP3").Value = [double]"0.2716857610474632"
$ws.Range("AQ3").Value = [double]"0.279179810725552"
$ws.Range("AR3").Value = [double]"0.2612312811980033"
$ws.Range("AS3").Value = [double]"0.2662337662337662"
$ws.Range("AT3").Value = [double]"0.2984822934232715"
$ws.Range("AU3").Value = [double]"0.2790697674418605"
$ws.Range("AV3").Value = [double]"0.2817589576547231"
$ws.Range("AW3").Value = [double]"0.2599337748344371"
$ws.Range("AX3").Value = [double]"0.284297520661157"
$ws.Range("AY3").Value = [double]"0.2669983416252073"
$ws.Range("AZ3").Value = [double]"0.281803167841809"
$ws.Range("BA3").Value = [double]"0.01236112433725989"

# Row 4
$ws.Range("A4").Value = [double]"2"
$ws.Range("B4").Value = [double]"0.02632591794319144"
$ws.Range("C4").Value = [double]"0.02756057908183015"
$ws.Range("D4").Value = [double]"0.02684049079754601"
$ws.Range("E4").Value = [double]"0.02829028290282903"
$ws.Range("F4").Value = [double]"0.02730899048787971"
$ws.Range("G4").Value = [double]"0.02780551501651432"
$ws.Range("H4").Value = [double]"0.02451731535396874"
$ws.Range("I4").Value = [double]"0.0271521135452021"
$ws.Range("J4").Value = [double]"0.02693915466790525"
$ws.Range("K4").Value = [double]"0.02443146896127843"
$ws.Range("L4").Value = [double]"0.02769066013921823"
$ws.Range("M4").Value = [double]"0.0259222333000997"
$ws.Range("N4").Value = [double]"0.02688664833486099"
$ws.Range("O4").Value = [double]"0.02921400417342917"
$ws.Range("P4").Value = [double]"0.02702910235736773"
$ws.Range("Q4").Value = [double]"0.0253418836436684"
$ws.Range("R4").Value = [double]"0.02783609361635954"
$ws.Range("S4").Value = [double]"0.02561316361378454"
$ws.Range("T4").Value = [double]"0.0256311871690584"
$ws.Range("U4").Value = [double]"0.02649513745309748"
$ws.Range("V4").Value = [double]"0.027386602349646"
$ws.Range("W4").Value = [double]"0.02711189400321667"
$ws.Range("X4").Value = [double]"0.02648493765972276"
$ws.Range("Y4").Value = [double]"0.02585939919479715"
$ws.Range("Z4").Value = [double]"0.02619485294117647"
$ws.Range("AA4").Value = [double]"0.02598808879263671"
$ws.Range("AB4").Value = [double]"0.02752924982794219"
$ws.Range("AC4").Value = [double]"0.0238819948489815"
$ws.Range("AD4").Value = [double]"0.02491349480968858"
$ws.Range("AE4").Value = [double]"0.02710027100271003"
$ws.Range("AF4").Value = [double]"0.0274627786777752"
$ws.Range("AG4").Value = [double]"0.02642763310408778"
$ws.Range("AH4").Value = [double]"0.0274830139705321"
$ws.Range("AI4").Value = [double]"0.02696078431372549"
$ws.Range("AJ4").Value = [double]"0.02558602727133446"
$ws.Range("AK4").Value = [double]"0.02518040841394135"
$ws.Range("AL4").Value = [double]"0.02421796165489405"
$ws.Range("AM4").Value = [double]"0.0271961678127173"
$ws.Range("AN4").Value = [double]"0.02561561332621788"
$ws.Range("AO4").Value = [double]"0.02669959101782545"
$ws.Range("AP4").Value = [double]"0.02568267966272143"
$ws.Range("AQ4").Value = [double]"0.02714932126696833"
$ws.Range("AR4").Value = [double]"0.0242041162414245"
$ws.Range("AS4").Value = [double]"0.02522300830513688"
$ws.Range("AT4").Value = [double]"0.02716390423572744"
$ws.Range("AU4").Value = [double]"0.02597804236895005"
$ws.Range("AV4").Value = [double]"0.02681754766702836"
$ws.Range("AW4").Value = [double]"0.02437888198757764"
$ws.Range("AX4").Value = [double]"0.02659656718725839"
$ws.Range("AY4").Value = [double]"0.02496124031007752"
$ws.Range("AZ4").Value = [double]"0.02636536033571057"
$ws.Range("BA4").Value = [double]"0.001171080329096825"

# Row 5
$ws.Range("A5").Value = [double]"3"
$ws.Range("B5").Value = [double]"0.07194412107101281"
$ws.Range("C5").Value = [double]"0.07194412107101281"
$ws.Range("D5").Value = [double]"0.07194412107101281"
$ws.Range("E5").Value = [double]"0.07194412107101281"
$ws.Range("F5").Value = [double]"0.07194412107101281"
$ws.Range("G5").Value = [double]"0.07194412107101281"
$ws.Range("H5").Value = [double]"0.07194412107101281"
$ws.Range("I5").Value = [double]"0.07194412107101281"
$ws.Range("J5").Value = [double]"0.07194412107101281"
$ws.Range("K5").Value = [double]"0.07194412107101281"
$ws.Range("L5").Value = [double]"0.07194412107101281"
$ws.Range("M5").Value = [double]"0.07194412107101281"
$ws.Range("N5").Value = [double]"0.07194412107101281"
$ws.Range("O5").Value = [double]"0.07194412107101281"
$ws.Range("P5").Value = [double]"0.07194412107101281"
$ws.Range("Q5").Value = [double]"0.07194412107101281"
$ws.Range("R5").Value = [double]"0.07194412107101281"
$ws.Range("S5").Value = [double]"0.07194412107101281"
$ws.Range("T5").Value = [double]"0.07194412107101281"
$ws.Range("U5").Value = [double]"0.07194412107101281"
$ws.Range("V5").Value = [double]"0.07194412107101281"
$ws.Range("W5").Value = [double]"0.07194412107101281"
$ws.Range("X5").Value = [double]"0.07194412107101281"
$ws.Range("Y5").Value = [double]"0.07194412107101281"
$ws.Range("Z5").Value = [double]"0.07194412107101281"
$ws.Range("AA5").Value = [double]"0.07194412107101281"
$ws.Range("AB5").Value = [double]"0.07194412107101281"
$ws.Range("AC5").Value = [double]"0.07194412107101281"
$ws.Range("AD5").Value = [double]"0.07194412107101281"
$ws.Range("AE5").Value = [double]"0.07194412107101281"
$ws.Range("AF5").Value = [double]"0.07194412107101281"
$ws.Range("AG5").Value = [double]"0.07194412107101281"
$ws.Range("AH5").Value = [double]"0.07194412107101281"
$ws.Range("AI5").Value = [double]"0.07194412107101281"
$ws.Range("AJ5").Value = [double]"0.07194412107101281"
$ws.Range("AK5").Value = [double]"0.07194412107101281"
$ws.Range("AL5").Value = [double]"0.07194412107101281"
$ws.Range("AM5").Value = [double]"0.07194412107101281"
$ws.Range("AN5").Value = [double]"0.07194412107101281"
$ws.Range("AO5").Value = [double]"0.07194412107101281"
$ws.Range("AP5").Value = [double]"0.07194412107101281"
$ws.Range("AQ5").Value = [double]"0.07194412107101281"
$ws.Range("AR5").Value = [double]"0.07194412107101281"
$ws.Range("AS5").Value = [double]"0.07194412107101281"
$ws.Range("AT5").Value = [double]"0.07194412107101281"
$ws.Range("AU5").Value = [double]"0.07194412107101281"
$ws.Range("AV5").Value = [double]"0.07194412107101281"
$ws.Range("AW5").Value = [double]"0.07194412107101281"
$ws.Range("AX5").Value = [double]"0.07194412107101281"
$ws.Range("AY5").Value = [double]"0.07194412107101281"
$ws.Range("AZ5").Value = [double]"0.07194412107101289"
$ws.Range("BA5").Value = [double]"8.326672684688674e-17"
